$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10 (weekly update): pushes existing rows 10-15 down to 11-16.
$ws.Rows(10).Insert()

# Populate the newly inserted row 10 with this week's data.
$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "Macroferia Regional de Talca"
$ws.Range("C10").Value = "Maule"
$ws.Range("D10").Value = 44680
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100104
$ws.Range("H10").Value = "Frutos de pepita"
$ws.Range("I10").Value = 100104001
$ws.Range("J10").Value = "Granada"
$ws.Range("K10").Value = "Wonderfull"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 15000
$ws.Range("Q10").Value = "$/caja 15 kilos granel"
$ws.Range("R10").Value = "Provincia de Limarí"
$ws.Range("S10").Value = 1000
$ws.Range("T10").Value = 15
